$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# BOM correction: row 10 (C12;C14) changes from the KEMET 4.7uF/10V/X5R
# (C0805C475K8PAC7210) part to the Vishay 4.7uF/25V/X7R (VJ126Y475KXXTW1BC)
# part, which also uses a bigger SMD1206 footprint instead of SMD0805.
# ---------------------------------------------------------------------------

# Remember the "Hyperlink" cell style so we can restore it after re-creating
# the F10 hyperlink (Hyperlinks.Add() below registers a fresh, functionally
# identical style record instead of reusing the existing one).
$hyperlinkStyle = $ws.Range("F9").Style

# Set right-to-left so freshly interned shared-string ids land in the same
# order as the authored file (G,F,E,D,C -> ids 95..99).
$ws.Range("G10").Value = "SMD1206"
$ws.Range("F10").Value = "77-VJ126Y475KXXTW1BC"
$ws.Range("E10").Value = "Vishay"
$ws.Range("D10").Value = "VJ126Y475KXXTW1BC"
$ws.Range("C10").Value = "4,7uF/25V/X7R"

# ---------------------------------------------------------------------------
# Hyperlinks: every RefMouser cell (column F) that has a product link keeps
# its link; only the one for row 10 needs to point at the new Vishay part.
# The COM shim here only supports wiping the whole collection (per-item /
# per-range Delete() does not scope correctly), so rebuild it from scratch,
# preserving the original relative order and appending the new row-10 link
# last (this mirrors exactly how the rIds shift in the target file).
# ---------------------------------------------------------------------------
$links = @(
  @{cell="F2";  url="https://pt.mouser.com/ProductDetail/Cypress-Semiconductor/CY7C68013A-56LTXI?qs=sGAEpiMZZMtv%252bwxsgy%2fhiMKDp0R93bWCFyLtDLcbUXY%3d"},
  @{cell="F4";  url="https://pt.mouser.com/ProductDetail/AVX/F0805B0R50FSTR?qs=sGAEpiMZZMtxU2g%2f1juGqcIlPyQcNmks%252ba9%2fntX%252bKIE%3d"},
  @{cell="F5";  url="https://pt.mouser.com/ProductDetail/STMicroelectronics/LD1117S33CTR?qs=sGAEpiMZZMsGz1a6aV8DcJsN6alkZ8%2fO7DPMKNVBfQo%3d"},
  @{cell="F6";  url="https://pt.mouser.com/ProductDetail/KEMET/C0805C104M5RAC7210?qs=sGAEpiMZZMs0AnBnWHyRQFCCI5cSbRT%2fI0VuIi6eznI%3d"},
  @{cell="F7";  url="https://pt.mouser.com/ProductDetail/IQD/LFXTAL066898Reel?qs=sGAEpiMZZMsBj6bBr9Q9aR%2fuGiDjvlISxfF6EyUKh0QY%252bn1mC20qZQ%3d%3d"},
  @{cell="F8";  url="https://pt.mouser.com/ProductDetail/KEMET/C0805C120J5GAC7210?qs=sGAEpiMZZMs0AnBnWHyRQOf5HOpVaXbhAZEuEA3l7uw%3d"},
  @{cell="F11"; url="https://pt.mouser.com/ProductDetail/Yageo/RT0805FRE0710KL?qs=sGAEpiMZZMu61qfTUdNhGxEjuuBLd0B41oWml1pU2QE%3d"},
  @{cell="F12"; url="https://pt.mouser.com/ProductDetail/Yageo/RT0805FRE07100KL?qs=sGAEpiMZZMu61qfTUdNhGxEjuuBLd0B4mFnwhs6X5Kw%3d"},
  @{cell="F9";  url="https://pt.mouser.com/ProductDetail/Yageo/RT0805DRD072K2L?qs=sGAEpiMZZMu61qfTUdNhGxEjuuBLd0B4GZ1bCnPt7TI%3d"},
  @{cell="F14"; url="https://pt.mouser.com/ProductDetail/Yageo/RC0805FR-070RL?qs=sGAEpiMZZMu61qfTUdNhG6gKAQVNBKOonL%252bE%2fLYSU34%3d"},
  @{cell="F15"; url="https://pt.mouser.com/ProductDetail/Kingbright/APT2012SGC?qs=sGAEpiMZZMvyj6n1w4pZD5QO%2fJ%252bHmnjZlXC14XYqTf8%3d"},
  @{cell="F16"; url="https://pt.mouser.com/ProductDetail/STMicroelectronics/M24128-DFMN6TP?qs=sGAEpiMZZMuVhdAcoizlRSnLhVfhSoFzqdBPoNcIiYg%3d"},
  @{cell="F18"; url="https://pt.mouser.com/ProductDetail/Texas-Instruments/SN74LVC245APWR?qs=sGAEpiMZZMs9F6aVvY09bnFMVNxjj%2fgIjJoUJoCdNnE%3d"},
  @{cell="F19"; url="https://pt.mouser.com/ProductDetail/STMicroelectronics/USBLC6-4SC6?qs=sGAEpiMZZMvxHShE6Whpu7uDi%252bhiPVjaQjZOzaiLtNs%3d"},
  @{cell="F3";  url="https://pt.mouser.com/ProductDetail/CUI/UJ2-MBH-1-SMT-TR?qs=sGAEpiMZZMulM8LPOQ%252byk%252br6FietFiXBYOK8b9%2fsuNIy0pbKqeLWLA%3d%3d"},
  @{cell="F20"; url="https://pt.mouser.com/ProductDetail/Harwin/M20-9740546?qs=sGAEpiMZZMs%252bGHln7q6pmzlZUuX%2f53qjhik29q3YN%2fs%3d"},
  @{cell="F10"; url="https://pt.mouser.com/ProductDetail/Vishay-Vitramon/VJ1206Y475KXXTW1BC?qs=sGAEpiMZZMs0AnBnWHyRQN7%2fAA2D2lPPu%252b7jaAAvyUeBMMU%2fzrdczQ%3d%3d"}
)

$tooltip = "Clique para visualizar informações adicionais sobre este produto"

$ws.Range("A1:G20").Hyperlinks.Delete()
foreach ($l in $links) {
  $ws.Hyperlinks.Add($ws.Range($l.cell), $l.url, [System.Reflection.Missing]::Value, $tooltip) | Out-Null
  # Hyperlinks.Add() registers a brand-new (if functionally identical)
  # "Hyperlink" style record on the target cell instead of reusing the
  # workbook's existing one; restore the original shared style so we don't
  # fork the style table.
  $ws.Range($l.cell).Style = $hyperlinkStyle
}

# ---------------------------------------------------------------------------
# Selection reflects the edited range in the saved view state.
# ---------------------------------------------------------------------------
$ws.Range("C10:G10").Select()
